$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "row1".."row4" labels in column G (rows 13-24) become plain numeric
# values (1-4) instead of text, so the now-unused shared strings disappear
# from sharedStrings.xml and every other shared-string index shifts down.
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("G16").Value = 4

$ws.Range("G17").Value = 1
$ws.Range("G18").Value = 2
$ws.Range("G19").Value = 3
$ws.Range("G20").Value = 4

$ws.Range("G21").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("G23").Value = 3
$ws.Range("G24").Value = 4

# The explanatory text in M14 is updated: "row=row1" -> "row=1"
$ws.Range("M14").Value = "select col1,col2,col3,col4 from tabella where user_id=tizio1 and row=1 mi restutisce tutta la prima riga di tizio1"

# Update the visible selection to match the author's saved view state.
[void]$ws.Range("N26").Select()
